$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date column C bumps from 46063 to 46064 for every data row
$ws.Range("C2").Value = 46064
$ws.Range("C3").Value = 46064
$ws.Range("C4").Value = 46064
$ws.Range("C5").Value = 46064
$ws.Range("C6").Value = 46064
$ws.Range("C7").Value = 46064
$ws.Range("C8").Value = 46064
$ws.Range("C9").Value = 46064

# Rows 4-9 get refreshed/reordered: two records (A 35734-2023 and A 5402-2026)
# drop to the bottom of the list while the remaining four shift up, keeping
# their relative order.

# Row 4 <- old row 5 data (A 25251-2025)
$ws.Range("A4").Value = "A 25251-2025"
$ws.Range("B4").Value = 45800.50082175926
$ws.Range("F4").Value = "Kommuner"
$ws.Range("G4").Value = 0.7

# Row 5 <- old row 7 data (A 6983-2023)
$ws.Range("A5").Value = "A 6983-2023"
$ws.Range("B5").Value = 44967.68585648148
$ws.Range("F5").ClearContents()
$ws.Range("G5").Value = 5.4

# Row 6 <- old row 8 data (A 25254-2025)
$ws.Range("A6").Value = "A 25254-2025"
$ws.Range("B6").Value = 45800.50479166667
$ws.Range("F6").Value = "Kommuner"
$ws.Range("G6").Value = 0.2

# Row 7 <- old row 9 data (A 26074-2025)
$ws.Range("A7").Value = "A 26074-2025"
$ws.Range("B7").Value = 45805.32366898148
$ws.Range("F7").ClearContents()
$ws.Range("G7").Value = 1.3

# Row 8 <- old row 4 data (A 35734-2023)
$ws.Range("A8").Value = "A 35734-2023"
$ws.Range("B8").Value = 45147.89258101852
$ws.Range("F8").ClearContents()
$ws.Range("G8").Value = 5.9

# Row 9 <- old row 6 data (A 5402-2026)
$ws.Range("A9").Value = "A 5402-2026"
$ws.Range("B9").Value = 46050.49721064815
$ws.Range("F9").Value = "Kommuner"
$ws.Range("G9").Value = 0.7
